# "CL and WO Screen" - add PRODUCT_CODE / PRODUCT_DESC / LOAN_AC_NUMBER columns
# to the "Target assignment in bulk" sheet, inserted right before the existing
# WO_AMOUNT / OS_AMOUNT / WO_DATE columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns at G:I (old G,H,I -> J,K,L), shifting formats/widths too.
$ws.Range("G1:I1").EntireColumn.Insert()

# --- Copy formatting from neighbouring cells onto the freshly inserted cells ---
# Header row (row 1) uses the same "text header" style as columns A,C,D,E,F (style id 3)
$ws.Range("C1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)

# Data rows: G (PRODUCT_CODE) and I (LOAN_AC_NUMBER) use the bordered "text" style (style id 5),
# the same one used by columns A/C
$ws.Range("C2").Copy()
$ws.Range("G2:G4").PasteSpecial(-4122)
$ws.Range("I2:I4").PasteSpecial(-4122)

# Data rows: H (PRODUCT_DESC) uses the plain bordered style (style id 4), same as columns B/D/F
$ws.Range("D2").Copy()
$ws.Range("H2:H4").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Column widths: new G:I should match the existing D:F width (29.140625 -> ~28.33 chars)
$ws.Range("G1:I4").ColumnWidth = 28.33

# --- Header values ---
$ws.Range("G1").Value = "PRODUCT_CODE"
$ws.Range("H1").Value = "PRODUCT_DESC"
$ws.Range("I1").Value = "LOAN_AC_NUMBER"

# --- Data values ---
$ws.Range("G2").Value = "0749"
$ws.Range("H2").Value = "EBL-TWO WHEELER"
$ws.Range("I2").Value = "3011380057880"

$ws.Range("G3").Value = "0749"
$ws.Range("H3").Value = "EBL-TWO WHEELER"
$ws.Range("I3").Value = "0031010025874"

$ws.Range("G4").Value = "0749"
$ws.Range("H4").Value = "EBL-TWO WHEELER"
$ws.Range("I4").Value = "0041010006371"

# --- View / selection state ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I10").Select()
